# ForestHillObserved.xlsx update:
#  - Insert a new EM sample row for ForestHill2023IrrigationFull (sim "15")
#    at 2024-03-25 (serial 45365), pushing all later rows down by one.
#  - Update the _xlnm._FilterDatabase defined name so it still spans the
#    whole data block after the insert.
#  - Append a new trailing row for ForestHill2023IrrigationPartial (sim "16")
#    on the same new sample date (45365), extending the met/clock range.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert new row 22 (shifts old rows 22-41 down to 23-42) -----------
$ws.Rows.Item(22).Insert()

# New row 22 data: ForestHill2023IrrigationFull sample on 45365.
$ws.Range("A22").Value = "ForestHill2023IrrigationFull"
$ws.Range("B22").Value = 45365
$ws.Range("K22").Value = 133.20000000000002
$ws.Range("L22").Value = 108.51666666666667

# --- 2. Append a brand-new last data row (43) -------------------------------
# ForestHill2023IrrigationPartial sample on the same new date, 45365.
$ws.Range("A43").Value = "ForestHill2023IrrigationPartial"
$ws.Range("B43").Value = 45365
$ws.Range("B43").NumberFormat = "d-mmm-yy"
$ws.Range("K43").Value = 84.583333333333343
$ws.Range("L43").Value = 48.733333333333341

# --- 3. Keep the hidden AutoFilter defined name in sync with the new range -
$wb.Names.Item("CottonObserved!_FilterDatabase").RefersTo = "=CottonObserved!`$A`$1:`$EQ`$2582"

# --- 4. Leave the selection on the last cell touched, like the author did --
$ws.Range("K43").Select()
